$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Fgf8"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07200033333333333
$ws.Range("H2").Value = 0.216001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.662797333333334
$ws.Range("N2").Value = 13.988392
$ws.Range("O2").Value = 0.7324994586787992
$ws.Range("P2").Value = 0.7324994586787993
$ws.Range("Q2").Value = 0.3357229622657778
$ws.Range("R2").Value = 3.021506660392
$ws.Range("S2").Value = 0.7324994586787992
$ws.Range("T2").Value = 0.7324994586787993

$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Fgf8"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07200033333333333
$ws.Range("H3").Value = 0.216001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6655859999999999
$ws.Range("N3").Value = 1.996758
$ws.Range("O3").Value = 0.1045598489170565
$ws.Range("P3").Value = 0.1045598489170565
$ws.Range("Q3").Value = 0.04792241386199999
$ws.Range("R3").Value = 0.4313017247579999
$ws.Range("S3").Value = 0.1045598489170565
$ws.Range("T3").Value = 0.1045598489170565

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fgf8"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07200033333333333
$ws.Range("H4").Value = 0.216001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7894166666666665
$ws.Range("N4").Value = 2.36825
$ws.Range("O4").Value = 0.1240129561007488
$ws.Range("P4").Value = 0.1240129561007488
$ws.Range("Q4").Value = 0.05683826313888888
$ws.Range("R4").Value = 0.5115443682499999
$ws.Range("S4").Value = 0.1240129561007488
$ws.Range("T4").Value = 0.1240129561007488

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf8"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07200033333333333
$ws.Range("H5").Value = 0.216001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08057833333333334
$ws.Range("N5").Value = 0.241735
$ws.Range("O5").Value = 0.01265840681643176
$ws.Range("P5").Value = 0.01265840681643176
$ws.Range("Q5").Value = 0.005801666859444445
$ws.Range("R5").Value = 0.052215001735
$ws.Range("S5").Value = 0.01265840681643176
$ws.Range("T5").Value = 0.01265840681643176

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf8"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.07200033333333333
$ws.Range("H6").Value = 0.216001
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16722
$ws.Range("N6").Value = 0.50166
$ws.Range("O6").Value = 0.02626932948696365
$ws.Range("P6").Value = 0.02626932948696365
$ws.Range("Q6").Value = 0.01203989574
$ws.Range("R6").Value = 0.10835906166
$ws.Range("S6").Value = 0.02626932948696365
$ws.Range("T6").Value = 0.02626932948696365

# Remove the now-obsolete row 7 (data consolidated into rows 2-6)
$ws.Rows("7").Delete()
